# Fruta / hortaliza, semanal
#
# Inserts a new weekly price-report row for "Poroto granado" (Feria
# Lagunitas de Puerto Montt) above the current row 19. Excel's native row
# insert shifts the existing rows 19-28 down to 20-29 (and grows the used
# range from A1:R28 to A1:R29), after which the new row's cells are filled
# in with the latest report's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19:28 down to 20:29, leaving a blank row 19 to populate.
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44596
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112030
$ws.Range("G19").Value = "Poroto granado"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 31000
$ws.Range("L19").Value = 31000
$ws.Range("M19").Value = 31000
$ws.Range("N19").Value = "`$/saco 25 kilos"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 1240
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
